$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 9998.5
$ws.Range("I96").Value = 9998
$ws.Range("K96").Value = 29994
$ws.Range("M96").Value = -28621

$ws.Range("H127").Value = 1382.8
$ws.Range("J127").Value = 1957.5
$ws.Range("L127").Value = 5872.5
$ws.Range("N127").Value = -15792.5

$ws.Range("H129").Value = 1031.1428
$ws.Range("I129").Value = 1040.5
$ws.Range("K129").Value = 3121.5
$ws.Range("M129").Value = 1878.5

$ws.Range("H132").Value = 2016.4359
$ws.Range("I132").Value = 1532.8
$ws.Range("J132").Value = 6248.25
$ws.Range("K132").Value = 4598.4
$ws.Range("L132").Value = 18744.75
$ws.Range("M132").Value = -2068.4
$ws.Range("N132").Value = -23804.75

$ws.Range("H137").Value = 327337.12
$ws.Range("I137").Value = 2145.8262
$ws.Range("J137").Value = 683499
$ws.Range("K137").Value = 6437.4786
$ws.Range("L137").Value = 2050497
$ws.Range("M137").Value = -3887.4786
$ws.Range("N137").Value = -2055597

$ws.Range("H138").Value = 3491.9531
$ws.Range("I138").Value = 2944.0557
$ws.Range("J138").Value = 3706.348
$ws.Range("K138").Value = 8832.167099999999
$ws.Range("L138").Value = 11119.044
$ws.Range("M138").Value = -3692.167099999999
$ws.Range("N138").Value = -21399.044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 600
$ws.Range("I16").Value = 600
$ws.Range("K16").Value = 600
$ws.Range("M16").Value = -313

$ws.Range("H19").Value = 669.3333
$ws.Range("J19").Value = 500
$ws.Range("L19").Value = 500
$ws.Range("N19").Value = -958

$ws.Range("H32").Value = 22739.389
$ws.Range("I32").Value = 12453.3125
$ws.Range("K32").Value = 12453.3125
$ws.Range("M32").Value = -12166.3125

$ws.Range("H45").Value = 7496.9414
$ws.Range("I45").Value = 9595.833000000001
$ws.Range("K45").Value = 9595.833000000001
$ws.Range("M45").Value = -9218.833000000001

$ws.Range("H97").Value = 475.13043
$ws.Range("I97").Value = 554.1053000000001
$ws.Range("K97").Value = 554.1053000000001
$ws.Range("M97").Value = -58.10530000000006

$ws.Range("H122").Value = 3012.8
$ws.Range("I122").Value = 2887.1
$ws.Range("J122").Value = 3264.2
$ws.Range("K122").Value = 8661.299999999999
$ws.Range("L122").Value = 9792.599999999999
$ws.Range("M122").Value = -6211.299999999999
$ws.Range("N122").Value = -14692.6

$ws.Range("H132").Value = 2175.8484
$ws.Range("I132").Value = 1853.4348
$ws.Range("J132").Value = 2917.4
$ws.Range("K132").Value = 5560.3044
$ws.Range("L132").Value = 8752.200000000001
$ws.Range("M132").Value = -3030.3044
$ws.Range("N132").Value = -13812.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 200789
$ws.Range("J141").Value = 200789
$ws.Range("L141").Value = 200789
$ws.Range("N141").Value = -211149

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4851.354
$ws.Range("I31").Value = 3134.2222
$ws.Range("J31").Value = 5508.9785
$ws.Range("K31").Value = 3134.2222
$ws.Range("L31").Value = 5508.9785
$ws.Range("M31").Value = -2839.2222
$ws.Range("N31").Value = -6098.9785

$ws.Range("H34").Value = 4851.354
$ws.Range("I34").Value = 3134.2222
$ws.Range("J34").Value = 5508.9785
$ws.Range("K34").Value = 3134.2222
$ws.Range("L34").Value = 5508.9785
$ws.Range("M34").Value = -2932.2222
$ws.Range("N34").Value = -5912.9785

$ws.Range("H94").Value = 15064.667
$ws.Range("I94").Value = 24366.4
$ws.Range("K94").Value = 24366.4
$ws.Range("M94").Value = -23915.4

$ws.Range("H99").Value = 6251962.5
$ws.Range("I99").Value = 1156
$ws.Range("J99").Value = 10419167
$ws.Range("K99").Value = 1156
$ws.Range("L99").Value = 10419167
$ws.Range("M99").Value = 342
$ws.Range("N99").Value = -10422163

$ws.Range("H126").Value = 6251962.5
$ws.Range("I126").Value = 1156
$ws.Range("J126").Value = 10419167
$ws.Range("K126").Value = 3468
$ws.Range("L126").Value = 31257501
$ws.Range("M126").Value = -998
$ws.Range("N126").Value = -31262441

$ws.Range("H134").Value = 2056.8223
$ws.Range("I134").Value = 1817.5428
$ws.Range("J134").Value = 2894.3
$ws.Range("K134").Value = 5452.6284
$ws.Range("L134").Value = 8682.900000000001
$ws.Range("M134").Value = -2917.6284
$ws.Range("N134").Value = -13752.9

$ws.Range("H141").Value = 191943.28
$ws.Range("J141").Value = 191943.28
$ws.Range("L141").Value = 191943.28
$ws.Range("N141").Value = -202303.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3040.8235
$ws.Range("J39").Value = 3832.8333
$ws.Range("L39").Value = 11498.4999
$ws.Range("N39").Value = -12086.4999

$ws.Range("H55").Value = 90914180
$ws.Range("J55").Value = 125006430
$ws.Range("L55").Value = 375019290
$ws.Range("N55").Value = -375019644

$ws.Range("H107").Value = 422.55554
$ws.Range("J107").Value = 475.66666
$ws.Range("L107").Value = 1426.99998
$ws.Range("N107").Value = -5266.999980000001

$ws.Range("H129").Value = 2725.5
$ws.Range("I129").Value = 1776.3334
$ws.Range("J129").Value = 3674.6667
$ws.Range("K129").Value = 5329.0002
$ws.Range("L129").Value = 11024.0001
$ws.Range("M129").Value = -329.0002000000004
$ws.Range("N129").Value = -21024.0001

$ws.Range("H131").Value = 31021.97
$ws.Range("I131").Value = 92482.73
$ws.Range("J131").Value = 2852.4583
$ws.Range("K131").Value = 277448.19
$ws.Range("L131").Value = 8557.374899999999
$ws.Range("M131").Value = -272408.19
$ws.Range("N131").Value = -18637.3749

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1802
$ws.Range("I102").Value = 1717.6
$ws.Range("K102").Value = 1717.6
$ws.Range("M102").Value = -95.59999999999991

$ws.Range("H132").Value = 2153.5144
$ws.Range("I132").Value = 1838.7916
$ws.Range("K132").Value = 5516.3748
$ws.Range("M132").Value = -2986.3748

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26022.523
$ws.Range("J7").Value = 13709.917
$ws.Range("L7").Value = 13709.917
$ws.Range("N7").Value = -13933.917

$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 2500
$ws.Range("K22").Value = 2500
$ws.Range("M22").Value = -2205

$ws.Range("H25").Value = 2253.5
$ws.Range("I25").Value = 7
$ws.Range("J25").Value = 4500
$ws.Range("K25").Value = 7
$ws.Range("L25").Value = 4500
$ws.Range("M25").Value = 223
$ws.Range("N25").Value = -4960

$ws.Range("H27").Value = 2500
$ws.Range("I27").Value = 2500
$ws.Range("K27").Value = 2500
$ws.Range("M27").Value = -2393

$ws.Range("H55").Value = 1355.1
$ws.Range("I55").Value = 556.6667
$ws.Range("K55").Value = 556.6667
$ws.Range("M55").Value = -383.6667

$ws.Range("H68").Value = 6665.8887
$ws.Range("I68").Value = 6499.375
$ws.Range("K68").Value = 6499.375
$ws.Range("M68").Value = -5750.375

$ws.Range("H71").Value = 6665.8887
$ws.Range("I71").Value = 6499.375
$ws.Range("K71").Value = 32496.875
$ws.Range("M71").Value = -28752.875

$ws.Range("H126").Value = 26022.523
$ws.Range("J126").Value = 13709.917
$ws.Range("L126").Value = 41129.751
$ws.Range("N126").Value = -46069.751

$ws.Range("H137").Value = 65749.75
$ws.Range("J137").Value = 72666.336
$ws.Range("L137").Value = 72666.336
$ws.Range("N137").Value = -82866.336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 55328
$ws.Range("J109").Value = 61660
$ws.Range("L109").Value = 61660
$ws.Range("N109").Value = -64434

$ws.Range("H126").Value = 1821.0682
$ws.Range("I126").Value = 1608.6061
$ws.Range("K126").Value = 4825.8183
$ws.Range("M126").Value = -2355.8183

$ws.Range("H132").Value = 913934.0600000001
$ws.Range("I132").Value = 10942
$ws.Range("J132").Value = 2900516.5
$ws.Range("K132").Value = 32826
$ws.Range("L132").Value = 8701549.5
$ws.Range("M132").Value = -30296
$ws.Range("N132").Value = -8706609.5
